$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws 2 4 "291.09"
Set-TextCell $ws 2 5 "-6.06%"
Set-TextCell $ws 2 7 "15"

# Row 3
Set-TextCell $ws 3 4 "39.87"
Set-TextCell $ws 3 5 "-2.76%"
Set-TextCell $ws 3 7 "15"

# Row 4
Set-TextCell $ws 4 4 "5.010"
Set-TextCell $ws 4 5 "-3.45%"
Set-TextCell $ws 4 7 "15"

# Row 5
Set-TextCell $ws 5 4 "0.07373"
Set-TextCell $ws 5 5 "-4.02%"
Set-TextCell $ws 5 7 "15"

# Row 6
Set-TextCell $ws 6 4 "4.282"
Set-TextCell $ws 6 5 "-0.32%"
Set-TextCell $ws 6 7 "15"

# Row 7
Set-TextCell $ws 7 4 "1.548"
Set-TextCell $ws 7 5 "-8.94%"
Set-TextCell $ws 7 7 "15"

# Row 8
Set-TextCell $ws 8 4 "0.9128"
Set-TextCell $ws 8 5 "-2.07%"
Set-TextCell $ws 8 7 "15"

# Row 9
Set-TextCell $ws 9 4 "0.1194"
Set-TextCell $ws 9 5 "-6.80%"
Set-TextCell $ws 9 7 "15"

# Row 10
Set-TextCell $ws 10 4 "0.1756"
Set-TextCell $ws 10 5 "-4.69%"
Set-TextCell $ws 10 7 "15"

# Row 11
Set-TextCell $ws 11 4 "0.08725"
Set-TextCell $ws 11 5 "-4.54%"
Set-TextCell $ws 11 7 "15"

# Row 12
Set-TextCell $ws 12 4 "0.04155"
Set-TextCell $ws 12 5 "-1.72%"
Set-TextCell $ws 12 7 "15"

# Row 13
Set-TextCell $ws 13 4 "0.1053"
Set-TextCell $ws 13 5 "0.25%"
Set-TextCell $ws 13 7 "15"

# Row 14
Set-TextCell $ws 14 4 "0.001274"
Set-TextCell $ws 14 5 "-0.84%"
Set-TextCell $ws 14 7 "15"

# Row 15
Set-TextCell $ws 15 4 "0.005906"
Set-TextCell $ws 15 5 "0.34%"
Set-TextCell $ws 15 7 "15"

# Row 16
Set-TextCell $ws 16 4 "3.397"
Set-TextCell $ws 16 5 "1.41%"
Set-TextCell $ws 16 7 "15"

# Row 17
Set-TextCell $ws 17 7 "15"

# Row 18
Set-TextCell $ws 18 7 "15"

# Row 19
Set-TextCell $ws 19 4 "7.552"
Set-TextCell $ws 19 5 "-0.64%"
Set-TextCell $ws 19 7 "15"

# Row 20
Set-TextCell $ws 20 5 "0.61%"
Set-TextCell $ws 20 7 "15"

# Row 21
Set-TextCell $ws 21 5 "5.92%"
Set-TextCell $ws 21 7 "15"

# Row 22
Set-TextCell $ws 22 4 "0.03838"
Set-TextCell $ws 22 5 "-4.49%"
Set-TextCell $ws 22 7 "15"

# Row 23
Set-TextCell $ws 23 4 "0.001269"
Set-TextCell $ws 23 5 "0.27%"
Set-TextCell $ws 23 7 "15"

# Row 24
Set-TextCell $ws 24 4 "0.003896"
Set-TextCell $ws 24 5 "-4.87%"
Set-TextCell $ws 24 7 "15"

# Row 25
Set-TextCell $ws 25 7 "15"

# Row 26
Set-TextCell $ws 26 4 "0.0003725"
Set-TextCell $ws 26 7 "15"

# Row 27
Set-TextCell $ws 27 7 "15"

# Row 28
Set-TextCell $ws 28 7 "15"

# Row 29
Set-TextCell $ws 29 7 "15"

# Row 30
Set-TextCell $ws 30 7 "15"

# Row 31
Set-TextCell $ws 31 7 "15"

# Row 32
Set-TextCell $ws 32 7 "15"

# Row 33
Set-TextCell $ws 33 7 "15"

# Row 34
Set-TextCell $ws 34 7 "15"

# Row 35
Set-TextCell $ws 35 7 "15"

# Row 36
Set-TextCell $ws 36 7 "15"

# Row 37
Set-TextCell $ws 37 7 "15"

# Row 38
Set-TextCell $ws 38 4 "0.02341"
Set-TextCell $ws 38 5 "-8.29%"
Set-TextCell $ws 38 7 "15"

# Row 39
Set-TextCell $ws 39 4 "0.05022"
Set-TextCell $ws 39 5 "-5.51%"
Set-TextCell $ws 39 7 "15"

# Row 40
Set-TextCell $ws 40 4 "0.007690"
Set-TextCell $ws 40 5 "-1.95%"
Set-TextCell $ws 40 7 "15"

# Row 41
Set-TextCell $ws 41 5 "132.07%"
Set-TextCell $ws 41 7 "15"

# Row 42
Set-TextCell $ws 42 4 "0.1273"
Set-TextCell $ws 42 5 "-2.90%"
Set-TextCell $ws 42 7 "15"

# Row 43
Set-TextCell $ws 43 4 "0.007365"
Set-TextCell $ws 43 5 "10.82%"
Set-TextCell $ws 43 7 "15"

# Row 44
Set-TextCell $ws 44 4 "0.006983"
Set-TextCell $ws 44 5 "-13.99%"
Set-TextCell $ws 44 7 "15"

# Row 45
Set-TextCell $ws 45 4 "0.3152"
Set-TextCell $ws 45 5 "2.09%"
Set-TextCell $ws 45 7 "15"

# Row 46
Set-TextCell $ws 46 4 "0.00006514"
Set-TextCell $ws 46 5 "-4.08%"
Set-TextCell $ws 46 7 "15"

# Row 47
Set-TextCell $ws 47 4 "0.00000000750"
Set-TextCell $ws 47 5 "-0.02%"
Set-TextCell $ws 47 7 "15"

# Row 48
Set-TextCell $ws 48 5 "11.83%"
Set-TextCell $ws 48 7 "15"

# Row 49
Set-TextCell $ws 49 4 "0.004203"
Set-TextCell $ws 49 5 "35.42%"
Set-TextCell $ws 49 7 "15"

# Row 50
Set-TextCell $ws 50 4 "0.00002101"
Set-TextCell $ws 50 5 "-0.02%"
Set-TextCell $ws 50 7 "15"

# Row 51
Set-TextCell $ws 51 4 "0.0002001"
Set-TextCell $ws 51 5 "-0.02%"
Set-TextCell $ws 51 7 "15"
